$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as exact text, avoiding Excel auto-number conversion
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '43.943.42'
Set-TextCell $ws.Range('D3') '2.355.79'
Set-TextCell $ws.Range('E3') '  -0.42%  '
Set-TextCell $ws.Range('E4') '  -0.01%  '
Set-TextCell $ws.Range('D5') '0.675'
Set-TextCell $ws.Range('E5') '  -2.59%  '
Set-TextCell $ws.Range('E6') '  -1.04%  '
Set-TextCell $ws.Range('D7') '73.60'
Set-TextCell $ws.Range('E7') '  -0.92%  '
Set-TextCell $ws.Range('E8') '  -0.03%  '
Set-TextCell $ws.Range('D9') '0.604'
Set-TextCell $ws.Range('E9') '  +0.54%  '
Set-TextCell $ws.Range('E10') '  -2.63%  '
Set-TextCell $ws.Range('D11') '58.93'
Set-TextCell $ws.Range('E11') '  +1.95%  '
Set-TextCell $ws.Range('D12') '33.92'
Set-TextCell $ws.Range('E12') '  +7.11%  '
Set-TextCell $ws.Range('D13') '7.35'
Set-TextCell $ws.Range('E13') '  -2.55%  '
Set-TextCell $ws.Range('E14') '  +0.03%  '
Set-TextCell $ws.Range('D15') '2.705.60'
Set-TextCell $ws.Range('E15') '  -0.49%  '
Set-TextCell $ws.Range('D16') '16.40'
Set-TextCell $ws.Range('E16') '  -3.76%  '
Set-TextCell $ws.Range('D17') '0.914'
Set-TextCell $ws.Range('E17') '  -0.52%  '
Set-TextCell $ws.Range('D18') '2.351.02'
Set-TextCell $ws.Range('E18') '  -0.08%  '
Set-TextCell $ws.Range('D19') '43.856.01'
Set-TextCell $ws.Range('E19') '  -0.85%  '
Set-TextCell $ws.Range('E20') '  -2.02%  '
Set-TextCell $ws.Range('E21') '  -0.10%  '
Set-TextCell $ws.Range('D22') '77.71'
Set-TextCell $ws.Range('E22') '  -1.22%  '
Set-TextCell $ws.Range('D23') '257.31'
Set-TextCell $ws.Range('E23') '  -0.02%  '
Set-TextCell $ws.Range('E24') '  +16.96%  '
Set-TextCell $ws.Range('E25') '  -0.10%  '
Set-TextCell $ws.Range('D26') '3.74'
Set-TextCell $ws.Range('E26') '  -0.65%  '
Set-TextCell $ws.Range('D27') '2.50'
Set-TextCell $ws.Range('E27') '  -2.69%  '
Set-TextCell $ws.Range('D28') '10.62'
Set-TextCell $ws.Range('E28') '  -1.85%  '
Set-TextCell $ws.Range('D29') '2.28'
Set-TextCell $ws.Range('E29') '  -1.57%  '
Set-TextCell $ws.Range('D30') '22.72'
Set-TextCell $ws.Range('E30') '  -0.07%  '
Set-TextCell $ws.Range('D31') '178.05'
Set-TextCell $ws.Range('E31') '  +1.58%  '
Set-TextCell $ws.Range('E32') '  -0.66%  '
Set-TextCell $ws.Range('D33') '0.137'
Set-TextCell $ws.Range('E33') '  +0.38%  '
Set-TextCell $ws.Range('E34') '  -0.31%  '
Set-TextCell $ws.Range('D35') '5.22'
Set-TextCell $ws.Range('E35') '  -3.74%  '
Set-TextCell $ws.Range('D36') '5.44'
Set-TextCell $ws.Range('E36') '  +0.44%  '
Set-TextCell $ws.Range('E37') '  -2.26%  '
Set-TextCell $ws.Range('D38') '2.40'
Set-TextCell $ws.Range('E38') '  -3.21%  '
Set-TextCell $ws.Range('D39') '6.39'
Set-TextCell $ws.Range('E39') '  -2.52%  '
Set-TextCell $ws.Range('E40') '  +1.16%  '
Set-TextCell $ws.Range('D41') '66.94'
Set-TextCell $ws.Range('E41') '  +25.32%  '
Set-TextCell $ws.Range('E42') '  +10.47%  '
Set-TextCell $ws.Range('E43') '  +13.85%  '
Set-TextCell $ws.Range('D44') '9.14'
Set-TextCell $ws.Range('E44') '  -0.02%  '
Set-TextCell $ws.Range('D45') '19.12'
Set-TextCell $ws.Range('E45') '  -0.54%  '
Set-TextCell $ws.Range('E46') '  +1.49%  '
Set-TextCell $ws.Range('E47') '  +0.28%  '
Set-TextCell $ws.Range('E48') '  -0.30%  '
Set-TextCell $ws.Range('E49') '  +0.04%  '
Set-TextCell $ws.Range('B50') 'Aave'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range('D50') '99.47'
Set-TextCell $ws.Range('E50') '  -1.81%  '
Set-TextCell $ws.Range('B51') 'ARBITRUM'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range('D51') '1.16'
Set-TextCell $ws.Range('E51') '  -2.11%  '
